{"js": "// 1) Remove the \"Deliverables\" table that sat right after the\n//    \"What are your deliverables?\" heading.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nfor (let i = tables.items.length - 1; i >= 0; i--) {\n  tables.items[i].delete();\n}\nawait context.sync();\n\n// 2) Append a new \"Real World Value and Interpretation\" section at the\n//    very end of the document body, after the existing bullet about the\n//    potential R Markdown website.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(Word.RangeLocation.after);\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n  '<w:bookmarkStart w:id=\"31\" w:name=\"real-world-value-and-interpretation\"/>' +\n  '<w:bookmarkEnd w:id=\"31\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">Real World Value and Interpretation</w:t></w:r>' +\n  \"</w:p>\" +\n  '<w:p><w:pPr><w:pStyle w:val=\"FirstParagraph\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">These methodologies can be applied in the real world work environment, and at my current employment in numerous ways, a few are listed below:</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">* Can I predict when raw materials are going to be at risk of backorder based off sales for product and demand of consumer?</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">* Can i identify which of my suppliers may be consistently a &quot;high risk&quot; when it comes to providing materials or products late?</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">* Is there a way to identify simply based off legacy product backorder information, which new products (based off inventory, lead time needed, source/suppliers) can be predicted to be at risk for a backorder issue?</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">* Additional real world applications</w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"Deliverables\" table that sat right after the\n#    \"What are your deliverables?\" heading.\nwhile ($d.Tables.Count -gt 0) {\n    $d.Tables.Item(1).Delete()\n}\n\n# 2) Append a new \"Real World Value and Interpretation\" section at the\n#    very end of the document body, after the existing bullet about the\n#    potential R Markdown website.\n#    NB: re-derive Paragraphs from a fresh $d.Content range instead of\n#    reusing $d.Paragraphs -- after the table delete above, the cached\n#    $d.Paragraphs collection can return stale/empty ranges.\n$lastPara = $d.Content.Paragraphs.Last\n$r = $lastPara.Range\n$r.InsertParagraphAfter()\n$newEmpty = $d.Content.Paragraphs.Last\n$target = $newEmpty.Range\n\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n'<w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr><w:bookmarkStart w:id=\"31\" w:name=\"real-world-value-and-interpretation\"/><w:bookmarkEnd w:id=\"31\"/><w:r><w:t xml:space=\"preserve\">Real World Value and Interpretation</w:t></w:r></w:p>' + `\n'<w:p><w:pPr><w:pStyle w:val=\"FirstParagraph\"/></w:pPr>' + `\n'<w:r><w:t xml:space=\"preserve\">These methodologies can be applied in the real world work environment, and at my current employment in numerous ways, a few are listed below:</w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\">* Can I predict when raw materials are going to be at risk of backorder based off sales for product and demand of consumer?</w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\">* Can i identify which of my suppliers may be consistently a &quot;high risk&quot; when it comes to providing materials or products late?</w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\">* Is there a way to identify simply based off legacy product backorder information, which new products (based off inventory, lead time needed, source/suppliers) can be predicted to be at risk for a backorder issue?</w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n'<w:r><w:t xml:space=\"preserve\">* Additional real world applications</w:t></w:r>' + `\n'</w:p>' + `\n'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($ooxml)\n\n# InsertXML inserts the new content before the pre-existing (now empty)\n# paragraph mark that InsertParagraphAfter created; remove that now-spare\n# trailing empty paragraph so the new section ends the document cleanly.\n$paras = $d.Content.Paragraphs\n$secondLast = $paras.Item($paras.Count - 1)\n$lastP = $paras.Last\n$delRange = $d.Range($secondLast.Range.End - 1, $lastP.Range.End)\n$delRange.Delete()\n"}
